$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D, shifting D:K (and beyond) to E:L.
$ws.Columns("D").Insert()

# Re-apply number formats to the freshly inserted column D cells so they
# match the rest of their row (date format for the three header rows,
# thousands format for every data row).
$ws.Range("D7,D38,D80").NumberFormat = "[$-409]d\-mmm\-yy;@"
$ws.Range("D8,D9,D10,D11,D12,D13,D14,D15,D16,D17,D18,D19,D20,D21,D22,D23,D24,D25,D26,D27,D28,D29,D30,D31,D32,D33,D34,D35,D39,D40,D41,D42,D43,D44,D45,D46,D47,D48,D49,D50,D51,D52,D53,D54,D55,D56,D57,D58,D59,D60,D61,D62,D63,D64,D65,D66,D67,D68,D69,D70,D71,D72,D73,D74,D75,D76,D77,D81,D82,D83,D84,D85,D86,D87,D88,D89,D90,D91,D92,D93,D94,D95,D96,D97,D98,D99,D100,D101,D102").NumberFormat = "#,##0"

# Set the new column D values (newest year of financial data).
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 142400
$ws.Range("D9").Value = 72500
$ws.Range("D10").Value = 69900
$ws.Range("D12").Value = 21900
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 300
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 136800
$ws.Range("D18").Value = 5600
$ws.Range("D20").Value = 2100
$ws.Range("D21").Value = 12700
$ws.Range("D22").Value = 700
$ws.Range("D23").Value = 7000
$ws.Range("D24").Value = -5400
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 12400
$ws.Range("D27").Value = 12400
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -2100
$ws.Range("D33").Value = 12400
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 12400
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 74100
$ws.Range("D42").Value = 9000
$ws.Range("D43").Value = 27600
$ws.Range("D44").Value = 30000
$ws.Range("D45").Value = "NA"
$ws.Range("D46").Value = 140800
$ws.Range("D47").Value = 44600
$ws.Range("D48").Value = 15000
$ws.Range("D49").Value = 6100
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 8400
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 214800
$ws.Range("D57").Value = 16600
$ws.Range("D58").Value = 0
$ws.Range("D59").Value = 16600
$ws.Range("D60").Value = 33200
$ws.Range("D61").Value = 0
$ws.Range("D62").Value = 2500
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 35700
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 22600
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 179100
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 12400
$ws.Range("D83").Value = 5000
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 33400
$ws.Range("D91").Value = -7300
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = 16700
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 5500
$ws.Range("D101").Value = -100
$ws.Range("D102").Value = 55500

# A handful of rows additionally had their prior-year figures restated
# (not just shifted) -- correct column E, and one deeper cell, to match.
$ws.Range("E20").Value = 1200
$ws.Range("E21").Value = 4000
$ws.Range("E22").Value = 800
$ws.Range("E32").Value = -1200
$ws.Range("E91").Value = -5700
$ws.Range("I91").Value = -1200
